$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "IPP369-Wind_1-ESS_1"
$ws.Range("C2").Value = 110.2263153116867
$ws.Range("D2").Value = 35.17716761216561
$ws.Range("E2").Value = 7912.102334764987
$ws.Range("F2").Value = 8912.102334764986
$ws.Range("G2").Value = 1473690210.906502
$ws.Range("H2").Value = 64.99999999999994
$ws.Range("I2").Value = 186257.7287999998
$ws.Range("J2").Value = 30.000000000019

# Row 3
$ws.Range("A3").Value = "IPP369-Wind_1-ESS_2"
$ws.Range("C3").Value = 95.99196815140722
$ws.Range("D3").Value = 53.15942122632188
$ws.Range("E3").Value = 9115.783038975145
$ws.Range("F3").Value = 10115.78303897515
$ws.Range("G3").Value = 1697885045.07307
$ws.Range("H3").Value = 64.99999999999991
$ws.Range("I3").Value = 186257.7287999997
$ws.Range("J3").Value = 30.00000000002
